# Rom Settings.xlsx — add a new "MultiBranchSupportUser" config row so the
# inventory transaction register resolves the management-role support
# contact instead of silently disappearing for that role.
#
# Before:
#   Row1: ID | Settings Name | Settings Value           (header)
#   Row2:    | domain_name   | http://rom_site:8000     (existing setting)
#
# After:
#   Row1: ID | Settings Name | Settings Value
#   Row2:    | domain_name            | http://rom_site:8000
#   Row3:    | MultiBranchSupportUser | manage@email.com   <-- new

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new settings row right under the existing "domain_name" row,
# mirroring its layout: column B holds the setting name, column C its value.
$ws.Range("B3").Value = "MultiBranchSupportUser"
$ws.Range("C3").Value = "manage@email.com"

# Match the formatting of the row immediately above so the new row renders
# identically to the other data rows (unbolded, default alignment/locking).
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false | Out-Null

# Leave the cursor where the source workbook's last save left it.
$ws.Range("C7").Select() | Out-Null
